$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Reset the "user answer" column (M) for the existing three questions back
#    to "NA" (the final-calculation column N recomputes to 99 automatically
#    via its existing IF formula).
# ---------------------------------------------------------------------------
$ws.Range("M2").Value = "NA"
$ws.Range("M3").Value = "NA"
$ws.Range("M4").Value = "NA"

# ---------------------------------------------------------------------------
# 2) Add two new dynamic questions (row 5 = slider question, row 6 = checkbox
#    question) under a new "סביבה" (environment) group.
#    Shared-string introduction order matters for a minimal diff, so the
#    "question type" (O column) cells are written first, then the question
#    text / group, then the checkbox answers.
# ---------------------------------------------------------------------------

# -- question-type column first --
$ws.Range("O5").Value = "סליידר"
$ws.Range("O6").Value = "check box"

# -- row 5: slider question about plant-based diet percentage --
$q5 = "מזון מהצומח כולל: דגנים, קטניות, פירות וירקות.  `nמזון מהחי כולל: בשר, עוף, דגים, ביצים, מוצרי חלב.  `nלפי ההגדרה הזאת, איזה אחוז מהתזונה שלך מבוסס על מזונות מהצומח? החליקי למענה"
$ws.Range("A5").Value = $q5
$ws.Range("B5").Value = "סביבה"

# -- row 6: checkbox question about where organic food is purchased --
$q6 = "בקניית מזון אורגני, אנא צייני איפה את רוכשת את המוצרים:`n(אפשר לסמן יותר מתשובה אחת)"
$ws.Range("A6").Value = $q6
$ws.Range("B6").Value = "סביבה"
$ws.Range("C6").Value = "ישירות מהחקלאי"
$ws.Range("D6").Value = "בחנות קטנה בעיר"
$ws.Range("E6").Value = "בעסק חברתי"
$ws.Range("F6").Value = "ברשת שיווק"
$ws.Range("G6").Value = "גידול עצמי"

# -- remaining "answer slot" / bookkeeping cells, all default to "NA" --
$ws.Range("C5:M5").Value = "NA"
$ws.Range("N5").Formula = "=M5"

$ws.Range("H6:N6").Value = "NA"

# -- "dependent question?" + trailing bookkeeping columns --
$ws.Range("P5").Value = "לא"
$ws.Range("Q5").Value = "NA"
$ws.Range("R5").Value = "NA"

$ws.Range("P6").Value = "לא"
$ws.Range("Q6").Value = "NA"

# -- wrap the long question text in col A, matching the other question rows --
$ws.Range("A5").WrapText = $true
$ws.Range("A6").WrapText = $true

# -- row heights sized for the wrapped question text --
$ws.Rows.Item(5).RowHeight = 126
$ws.Rows.Item(6).RowHeight = 78.75

# ---------------------------------------------------------------------------
# 3) Add a dropdown list data validation on column O ("סוג שאלה" / question
#    type) covering the whole column, offering the three known question
#    types.
# ---------------------------------------------------------------------------
$dvRange = $ws.Range("O2:O1048576")
$dvRange.Validation.Add(3, 1, 1, '"אמריקאית,סליידר,check box"')
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 4) Leave the final selection on Q6, matching where editing ended.
# ---------------------------------------------------------------------------
$ws.Range("Q6").Select()
